# Edit LOM3058.xlsx worksheet to match target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the "1033242 - Fabio Herbst Florenzano" professor-name
# value (B13/C13, no label in A13) is removed; rows below shift up by one.
$ws.Rows(13).Delete()

# After the shift, several B/C (value) cells need new text so each row's
# label (column A) and its value are re-paired per the target revision.
$ws.Range("B10").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C10").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"

$ws.Range("B18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C18").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."
$ws.Range("C19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."

$ws.Range("B20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."
$ws.Range("C20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."

$ws.Range("B21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
$ws.Range("C21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."

